# 임시 [AutoInput] 코드 작성
# Rename several "완료 요청" (complete REQUEST) labels to "완료 보고" (complete REPORT),
# add CellID/TrayID qualifiers to several branch-1/branch-2 automatic/manual cell-input
# steps, and assign short "35_1_N" step codes to the branch-1 "자동 투입" (auto input) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DSF Sequence")

# ---- "완료 요청" -> "완료 보고" wording fixes ----
$ws.Range("D18").Value = "Rack 출고 완료 보고 with TrayId/RackId"
$ws.Range("D19").Value = "Rack 출고 완료 보고 Confirm"
$ws.Range("D20").Value = "투입 Tray 투입 완료 보고"
$ws.Range("D21").Value = "투입 Tray 투입 완료 보고 Confirm"
$ws.Range("D28").Value = "배출 Tray 투입 완료 보고"
$ws.Range("D29").Value = "배출 Tray 투입 완료 보고 Confirm"
$ws.Range("D34").Value = "NG Tray Break 완료 보고 with TrayId"
$ws.Range("D35").Value = "NG Tray Break 완료 보고 Confirm"
$ws.Range("D36").Value = "NG Tray 투입 완료 보고"
$ws.Range("D37").Value = "NG Tray 투입 완료 보고 Confirm"

# ---- 분기1 (branch 1 - auto input) rows: add step codes in column B, and append
#      "with CellID" / "& TrayID" qualifiers to the relevant event descriptions ----
$ws.Range("B38").Value = "35_1_1"
$ws.Range("D38").Value = "[분기1 - 자동 투입] 자동 셀 투입 요청 with CellID & TrayID"

$ws.Range("B39").Value = "35_1_2"

$ws.Range("B40").Value = "35_1_3"
$ws.Range("D40").Value = "[분기1 - 자동 투입] 투입셀 정보 검사 요청 with CellID & TrayID"

$ws.Range("B41").Value = "35_1_4"

$ws.Range("B42").Value = "35_1_5"
$ws.Range("D42").Value = "[분기1 - 자동 투입] 셀 공정 시작 보고 with CellID"

$ws.Range("B43").Value = "35_1_6"

$ws.Range("B44").Value = "35_1_7"
$ws.Range("D44").Value = "[분기1 - 자동 투입] 자동 셀 투입 완료보고 with CellID & TrayID"

$ws.Range("B45").Value = "35_1_8"

# ---- 분기2 (branch 2 - manual input) : add "& TrayID" qualifier ----
$ws.Range("D48").Value = "[분기2 - 수동 투입] 투입셀 정보 검사 요청 with CellID & TrayID"

# ---- Restore the view/selection the author ended the edit session on ----
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D19").Select()
